# Generate Report for Handoff
#
# A new handoff was generated for the "9b2ca098-..." localization module.
# This refreshes the "Latest Handoff Datetime" cell (column D, row 6) on
# both the "zh-cn" and "de-de" status sheets to the newer handoff
# timestamps, while every other reported value for that run stays as-is.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D6").Value = "2016-03-10 18:38:52"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D6").Value = "2016-03-10 18:38:57"
